$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Thông Tin")
$ws2 = $wb.Worksheets.Item("Chức danh")

# ---------------------------------------------------------------------------
# Sheet "Thông Tin": fill in row 4 (previously partially blank) and build out
# rows 5-18 as copies of the templated employee row.
# ---------------------------------------------------------------------------

# --- Row 4: fill previously-blank cells, keep existing styling -------------
$ws1.Range("B4").Value = "Lê Thành Khang'"
$ws1.Range("C4").Value = "070288372732"
$ws1.Range("F4").Value = "0367909248"

# H4 needs a mailto hyperlink like H3. Adding a hyperlink resets the cell
# style, so stash/restore the original formatting around the call.
$ws1.Range("H4").Copy()
$ws1.Range("Z4").PasteSpecial(-4122)
$ws1.Range("H4").Value = "tuan.197pm21996@vanlanguni.vn"
$ws1.Hyperlinks.Add($ws1.Range("H4"), "mailto:tuan.197pm21996@vanlanguni.vn", "", "", "tuan.197pm21996@vanlanguni.vn")
$ws1.Range("Z4").Copy()
$ws1.Range("H4").PasteSpecial(-4122)
$ws1.Range("Z4").Clear()

# --- Row 5: fix the one-off "last row" style on A5, fill blank cells -------
$ws1.Range("A3").Copy()
$ws1.Range("A5").PasteSpecial(-4122)

$ws1.Range("B5").Value = "Lê Thành Khang'"
$ws1.Range("C5").Value = "070288372732"
$ws1.Range("F5").Value = "0367909248"

$ws1.Range("H5").Copy()
$ws1.Range("Z5").PasteSpecial(-4122)
$ws1.Range("H5").Value = "tuan.197pm21996@vanlanguni.vn"
$ws1.Hyperlinks.Add($ws1.Range("H5"), "mailto:tuan.197pm21996@vanlanguni.vn", "", "", "tuan.197pm21996@vanlanguni.vn")
$ws1.Range("Z5").Copy()
$ws1.Range("H5").PasteSpecial(-4122)
$ws1.Range("Z5").Clear()

# --- Rows 6-18: duplicate row 5 (values + formatting) down the table -------
for ($r = 6; $r -le 18; $r++) {
    $ws1.Range("A5:M5").Copy($ws1.Range("A$r`:M$r"))
}

# --- Fix the STT (row number) column and re-add hyperlinks on rows 6-18 ----
for ($r = 4; $r -le 18; $r++) {
    $ws1.Range("A$r").Value = $r - 2
}

for ($r = 6; $r -le 18; $r++) {
    $ws1.Range("H$r").Copy()
    $ws1.Range("Z$r").PasteSpecial(-4122)
    $ws1.Range("H$r").Value = "tuan.197pm21996@vanlanguni.vn"
    $ws1.Hyperlinks.Add($ws1.Range("H$r"), "mailto:tuan.197pm21996@vanlanguni.vn", "", "", "tuan.197pm21996@vanlanguni.vn")
    $ws1.Range("Z$r").Copy()
    $ws1.Range("H$r").PasteSpecial(-4122)
    $ws1.Range("Z$r").Clear()
}

# --- Expand the worksheet Table (Table4) to cover the new rows -------------
$ws1.ListObjects.Item("Table4").Resize($ws1.Range("A2:M18"))

# ---------------------------------------------------------------------------
# Sheet "Chức danh": rename the job-title entry.
# ---------------------------------------------------------------------------
$ws2.Range("B5").Value = "Full-Stack Developer"

Write-Host "done"
